# Path To Graduation 2 - update the CPSC/CYBR plan to a Kinesiology-track
# plan (POLS/KINS/DSCI/CPSC courses) and extend the sheet with the
# Fall 2024 / Spring 2024 / Summer 2024 and Fall 2025 / Spring 2025 /
# Summer 2025 semester blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Fall 2022 block (rows 4-10), replacing the old CPSC/CYBR courses ----
$ws.Range("A4").Value = "POLS 1101"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "KINS 2271"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "CPSC 3121"
$ws.Range("F4").Value = 3

$ws.Range("A5").Value = "KINS 1105"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "KINS 2272"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "KINS 3126"
$ws.Range("F5").Value = 2

$ws.Range("A6").Value = "KINS 1106"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "KINS 3105"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "KINS 3165"
$ws.Range("F6").Value = 2

$ws.Range("A7").Value = "KINS 2105"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "KINS 3107"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "KINS 3256"
$ws.Range("F7").Value = 2

$ws.Range("A8").Value = "KINS 2135"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "DSCI 3111"
$ws.Range("D8").Value = 3

$ws.Range("A9").Value = "KINS 2345"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "CPSC 4000"
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = "KINS 2379"
$ws.Range("B10").Value = 1

# Row 11 totals already hold =SUM(B4:B10)/=SUM(D4:D10)/=SUM(F4:F10) - untouched.

# ---- Fall 2023 block (rows 13-19), replacing the old CPSC/CYBR courses ----
$ws.Range("A13").Value = "KINS 3127"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "KINS 3257"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "CPSC 4148"
$ws.Range("F13").Value = 3

$ws.Range("A14").Value = "CPSC 3165"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "KINS 3258"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "CPSC 4698"
$ws.Range("F14").Value = 3

$ws.Range("A15").Value = "KINS 3218"
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "KINS 3365"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "CPSC 4899"
$ws.Range("F15").Value = 3

$ws.Range("A16").Value = "KINS 3235"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "CPSC 3415"
$ws.Range("D16").Value = 1

$ws.Range("A17").Value = "KINS 3255"
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "CYBR 4125"
$ws.Range("D17").Value = 3

$ws.Range("A18").Value = "KINS 3316"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "CPSC 4135"
$ws.Range("D18").Value = 3

$ws.Range("C19").Value = "CYBR 4145"
$ws.Range("D19").Value = 3

# Row 20 totals already hold =SUM(B13:B19)/=SUM(D13:D19)/=SUM(F13:F19) - untouched.

# ---- New Fall 2024 / Spring 2024 / Summer 2024 block (rows 21-29) ----
$ws.Range("A21").Value = "Fall 2024"
$ws.Range("B21").Value = "Credits"
$ws.Range("C21").Value = "Spring 2024"
$ws.Range("D21").Value = "Credits"
$ws.Range("E21").Value = "Summer 2024"
$ws.Range("F21").Value = "Credits"

$ws.Range("A22").Value = "CPSC 4155"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "CPSC 4176"
$ws.Range("D22").Value = 3

$ws.Range("A23").Value = "CPSC 4157"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = "CPSC 4555"
$ws.Range("D23").Value = 3

$ws.Range("A24").Value = "CPSC 4175"
$ws.Range("B24").Value = 3

$ws.Range("A25").Value = "CPSC 4205"
$ws.Range("B25").Value = 3

$ws.Range("A26").Value = "CYBR 4416"
$ws.Range("B26").Value = 1

$ws.Range("A29").Value = "Total"
$ws.Range("B29").Formula = "=SUM(B22:B28)"
$ws.Range("C29").Value = "Total"
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Value = "Total"
$ws.Range("F29").Formula = "=SUM(F22:F28)"

# ---- New Fall 2025 / Spring 2025 / Summer 2025 block (rows 30-38) ----
$ws.Range("A30").Value = "Fall 2025"
$ws.Range("B30").Value = "Credits"
$ws.Range("C30").Value = "Spring 2025"
$ws.Range("D30").Value = "Credits"
$ws.Range("E30").Value = "Summer 2025"
$ws.Range("F30").Value = "Credits"

$ws.Range("A38").Value = "Total"
$ws.Range("B38").Formula = "=SUM(B31:B37)"
$ws.Range("C38").Value = "Total"
$ws.Range("D38").Formula = "=SUM(D31:D37)"
$ws.Range("E38").Value = "Total"
$ws.Range("F38").Formula = "=SUM(F31:F37)"
